# Added team record (Wins/Losses/Ties) columns to the Oakland Athletics
# 2018 roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting of the last existing header cell (AC1) onto
# the three new header cells so they pick up the same bold/centered/bordered
# style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player on the roster shares the team's overall 2018 record.
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 97   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 65   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
